$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix capitalization: "git" -> "Git" in the reflection text of the first diary entry (row 10)
$ws.Range("F10").Value = "Having little experience with Git itself, navigating the system was quite confusing; however, as discussed in the lecture, this seems like a necessary skill that requires plenty of practice to get used to"

# Revert the second diary entry (row 11) back to the blank template placeholder text.
# A11/B11 previously held real Date/Time values; switch their format back to Text so
# they can hold the placeholder strings instead (matches the styling used by the
# still-blank template rows below).
$ws.Range("A11:B11").NumberFormat = "@"

$ws.Range("A11").Value = "<what day?>"
$ws.Range("B11").Value = "<what time?>"
$ws.Range("C11").Value = "<as applicable, with whom?>"
$ws.Range("D11").Value = "<what did you want to accomplish?>"
$ws.Range("E11").Value = "<what did you actually accomplish?>"
$ws.Range("F11").Value = "<what insight(s) did you gain?>"
$ws.Range("G11").Value = "<how did you feel during the activity?>"

# Duplicate the same placeholder text into row 12 (the next blank template row)
$ws.Range("A12").Value = "<what day?>"
$ws.Range("B12").Value = "<what time?>"
$ws.Range("C12").Value = "<as applicable, with whom?>"
$ws.Range("D12").Value = "<what did you want to accomplish?>"
$ws.Range("E12").Value = "<what did you actually accomplish?>"
$ws.Range("F12").Value = "<what insight(s) did you gain?>"
$ws.Range("G12").Value = "<how did you feel during the activity?>"

# Row 13 col A gets the trailing "Etc." marker
$ws.Range("A13").Value = "Etc."
